$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("metadata")

# Add "Do not edit this sheet." text with yellow fill across E1:G1
$ws2.Range("E1").Value = "Do not edit this sheet."
$ws2.Range("E1:G1").Interior.Color = 65535

# Make metadata the active sheet/tab and select G9
$ws2.Activate()
[void]$ws2.Range("G9").Select()
